# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

# Map of row -> new value for sheet "展览"
$exhibitionUpdates = @{
    2  = 137
    4  = 431
    7  = 2185
    11 = 4941
    14 = 307
    17 = 186
    21 = 3899
    22 = 713
    23 = 674
    27 = 121
    31 = 580
    34 = 961
    35 = 2471
}

# Map of row -> new value for sheet "全部类型"
$allTypesUpdates = @{
    2  = 137
    4  = 431
    7  = 2185
    11 = 4941
    14 = 307
    17 = 186
    21 = 3899
    22 = 713
    23 = 674
    27 = 121
    31 = 580
    35 = 961
    36 = 2471
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
